# clientes.xlsx - "solucion de interconexion que presentaba problemas"
#
# The sheet stores a flat list of "cliente" records. Each client occupies a
# pair of rows: the first row holds the client name (col A) and description
# (col B) plus the first few (id, price) pairs starting at column C; the
# second row continues the (id, price) pairs if there are more than fit on
# the first row. Row 1 is special - it is the "default" client whose (id,
# price) list was extended with 6 more pairs (ids 7-12).
#
# This script applies the following fixes:
#  1. Row 1 (Clientes Varios / default client): append 6 more (id, price)
#     pairs in columns Q:AB.
#  2. Row 4: drop the stray trailing (id=-100) sentinel pair that had no
#     matching price - interconnection glitch duplicate.
#  3. Row 7 (Leonardo): extend the client description to mention "y
#     ponedora".
#  4. Row 8: this continuation row was empty; it now carries a price value
#     that belongs to the Leonardo record.
#  5. Row 25 (Rigoberto Padilla): the (id, price) pairs had drifted out of
#     sync across the interconnection - re-sequence them so the pair set is
#     complete (2/140, 12/188).
#  6. Add a new example client in rows 27-28 (Cliente de ejemplo /
#     descripcion de ejemplo).
#  7. Add a new client "Mario Leo" in rows 29-30 with his (id, price) pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 1: extend the default client's price table (ids 7-12) -------
$ws.Range("Q1").Value = 7
$ws.Range("R1").Value = 166
$ws.Range("S1").Value = 8
$ws.Range("T1").Value = 230
$ws.Range("U1").Value = 9
$ws.Range("V1").Value = 255
$ws.Range("W1").Value = 10
$ws.Range("X1").Value = 235
$ws.Range("Y1").Value = 11
$ws.Range("Z1").Value = 210
$ws.Range("AA1").Value = 12
$ws.Range("AB1").Value = 190

# --- 2. Row 4: remove the stray orphaned sentinel value ------------------
$ws.Range("I4").ClearContents()

# --- 3. Row 7: update Leonardo's description ------------------------------
$ws.Range("B7").Value = "Mayorista de`nalimento para `ngallinas y ponedora"

# --- 4. Row 8: continuation price for Leonardo ----------------------------
$ws.Range("C8").Value = 20

# --- 5. Row 25: re-sequence Rigoberto Padilla's (id, price) pairs --------
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 140
$ws.Range("E25").Value = 12
$ws.Range("F25").Value = 188

# --- 6. Rows 27-28: new example client ------------------------------------
$ws.Range("A27").Value = "Cliente de ejemplo"
$ws.Range("B27").Value = "descripcion de ejemplo"
$ws.Range("C28").Value = -100

# --- 7. Rows 29-30: new client "Mario Leo" --------------------------------
$ws.Range("A29").Value = "Mario Leo"
$ws.Range("B29").Value = "Cliente mayorista de lactomayma `ny lactomayma 22 porciento"
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 195
$ws.Range("E29").Value = 12
$ws.Range("F29").Value = 198
$ws.Range("C30").Value = 19
$ws.Range("D30").Value = 21
